$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per diff
$ws.Range("B2").Value = 5.8
$ws.Range("C2").Value = 10.7
$ws.Range("B3").Value = 5.7
$ws.Range("C3").Value = 9.3000000000000007

# Update the active selection to C3 (matches diff's <selection activeCell="C3" sqref="C3"/>)
$ws.Range("C3").Select()
